$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# New row for "The Magic of Thinking Big" goes right after the last existing row (90)
$lastRow = 90
$newRow = $lastRow + 1

# Copy the formatting of the previous row down into the new row first so the
# new cells pick up the same styles (e.g. date number format) without
# introducing any new style/numFmt entries.
$ws.Range("A" + $lastRow + ":G" + $lastRow).Copy() | Out-Null
$ws.Range("A" + $newRow + ":G" + $newRow).PasteSpecial(-4122) | Out-Null

$ws.Range("A" + $newRow).Value = "The Magic of Thinking Big"
$ws.Range("B" + $newRow).Value = "David Schwartz"
$ws.Range("C" + $newRow).Value = 43998
$ws.Range("D" + $newRow).Value = 44000
$ws.Range("E" + $newRow).Value = "confidence;self-improvement;goals;positivity"
$ws.Range("F" + $newRow).Value = "Audio"
$ws.Range("G" + $newRow).Value = "9 Hours 31 Mins"

# Match the scrolled/selected view state left behind after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
$ws.Range("C92").Select() | Out-Null
